$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reposition/resize the workbook window to match the latest save
$wbWin = $wb.Windows.Item(1)
$wbWin.Left = 1360
$wbWin.Top = 1840
$wbWin.Width = 24780
$wbWin.Height = 14600

# Append new rows of test data (18-21)
$ws.Range("A18").Value = 68
$ws.Range("B18").Value = 15
$ws.Range("C18").Value = 0
$ws.Range("D18").Value = 0

$ws.Range("A19").Value = 85
$ws.Range("B19").Value = 8
$ws.Range("C19").Value = 0
$ws.Range("D19").Value = 0

$ws.Range("A20").Value = 65
$ws.Range("B20").Value = 3
$ws.Range("C20").Value = 1
$ws.Range("D20").Value = 0

$ws.Range("A21").Value = 86
$ws.Range("B21").Value = 8
$ws.Range("C21").Value = 0
$ws.Range("D21").Value = 0
$ws.Range("E21").Value = "one cell near top edge was segmented smaller than in reality"

# Update freeze panes / scroll position and selection
$win = $excel.ActiveWindow
$win.FreezePanes = $false
$ws.Range("A11").Select() | Out-Null
$win.FreezePanes = $true
$ws.Range("C22").Select() | Out-Null
